$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 3197
$ws.Range("F7").Value = 322
$ws.Range("F8").Value = 7485
$ws.Range("F11").Value = 1216
$ws.Range("F14").Value = 605
$ws.Range("F15").Value = 1052
$ws.Range("F18").Value = 1217
$ws.Range("F20").Value = 5879
$ws.Range("F21").Value = 2327
$ws.Range("F22").Value = 4067
$ws.Range("F23").Value = 2264
$ws.Range("F24").Value = 226
$ws.Range("F29").Value = 66
$ws.Range("F34").Value = 558
$ws.Range("F35").Value = 332
$ws.Range("F36").Value = 262
$ws.Range("F37").Value = 947
$ws.Range("F39").Value = 73
$ws.Range("F41").Value = 246

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 578
$ws.Range("F11").Value = 61
$ws.Range("F14").Value = 127
$ws.Range("F16").Value = 1
$ws.Range("F18").Value = 124
$ws.Range("F25").Value = 3479
$ws.Range("F26").Value = 3479
$ws.Range("F31").Value = 29

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 1910
$ws.Range("F7").Value = 1894
$ws.Range("F8").Value = 2976
$ws.Range("F9").Value = 1177
$ws.Range("F10").Value = 1204
$ws.Range("F13").Value = 1944
$ws.Range("F14").Value = 8519
$ws.Range("F15").Value = 716

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 3197
$ws.Range("F5").Value = 1910
$ws.Range("F6").Value = 322
$ws.Range("F7").Value = 1178
$ws.Range("F8").Value = 1204
$ws.Range("F11").Value = 1216
$ws.Range("F14").Value = 716
$ws.Range("F16").Value = 578
$ws.Range("F17").Value = 578
$ws.Range("F19").Value = 605
$ws.Range("F20").Value = 1052
$ws.Range("F21").Value = 61
$ws.Range("F25").Value = 127
$ws.Range("F28").Value = 5879
$ws.Range("F29").Value = 2327
$ws.Range("F30").Value = 4067
$ws.Range("F31").Value = 2264
$ws.Range("F34").Value = 66
$ws.Range("F37").Value = 124
$ws.Range("F39").Value = 558
$ws.Range("F40").Value = 332
$ws.Range("F41").Value = 262
$ws.Range("F45").Value = 73
$ws.Range("F47").Value = 3479
$ws.Range("F50").Value = 29
